# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# Updates the FHIR StructureDefinition export: version bump, new publish date,
# publisher/jurisdiction metadata, and refreshed short/definition text for the
# root Extension element.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refresh publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# The old duplicate "Contact" / "No display for ContactDetail" row (row 10)
# becomes a "Jurisdiction" / "United States of America" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The second duplicate "Contact" row (row 11) is removed entirely.
$meta.Rows.Item(11).Delete()

# Elements sheet: update Short/Definition for the root Extension element (row 2)
$elements.Range("K2").Value = "Problem Priority"
$elements.Range("L2").Value = "Problem priority code"
